# "saved by Ebrahem Ayman"
# Adds a second submission row (name / email / repo link) to the
# "Open Source task" tracking sheet, with the email turned into a
# mailto: hyperlink (Excel's default "Hyperlink" cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$name  = "ابراهيم ايمن ابراهيم حسين"
$email = "ebrahemayman48667@gmail.com"
$repo  = "https://github.com/omargalal255/oos.git"

$ws.Range("A2").Value = $name
$ws.Range("B2").Value = $email
$ws.Range("C2").Value = $repo

# Turn the email cell into a clickable mailto: hyperlink - this is what
# mints the "Hyperlink" font/style in styles.xml and the <hyperlinks> part.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:$email")

# Move the active selection, matching the saved view state.
[void]$ws.Range("C5").Select()
